$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.135.60"
$ws.Range("E2").Value = "  -1.35%  "

$ws.Range("D3").Value = "1.794.40"
$ws.Range("E3").Value = "  -1.56%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'314.37"
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("D7").Value = "'0.5221"
$ws.Range("E7").Value = "  +2.20%  "

$ws.Range("D8").Value = "'0.3821"
$ws.Range("E8").Value = "  -3.18%  "

$ws.Range("D9").Value = "'0.07964"
$ws.Range("E9").Value = "  -3.00%  "

$ws.Range("D10").Value = "'41.41"
$ws.Range("E10").Value = "  -0.58%  "

$ws.Range("D11").Value = "'1.093"
$ws.Range("E11").Value = "  -1.57%  "

$ws.Range("D12").Value = "'6.279"
$ws.Range("E12").Value = "  -0.73%  "

$ws.Range("E13").Value = "  +0.16%  "

$ws.Range("D14").Value = "'20.49"
$ws.Range("E14").Value = "  -2.96%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.297"
$ws.Range("E15").Value = "  -2.90%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.796.80"
$ws.Range("E16").Value = "  -1.50%  "

$ws.Range("D17").Value = "'91.86"
$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("D18").Value = "'0.00001087"
$ws.Range("E18").Value = "  -3.67%  "

$ws.Range("D19").Value = "'0.06567"
$ws.Range("E19").Value = "  -1.36%  "

$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").Value = "'17.31"
$ws.Range("E21").Value = "  -2.75%  "

$ws.Range("D22").Value = "'5.945"
$ws.Range("E22").Value = "  -2.38%  "

$ws.Range("D23").Value = "28.177.64"
$ws.Range("E23").Value = "  -1.30%  "

$ws.Range("E24").Value = "  -2.05%  "

$ws.Range("D25").Value = "'2.264"
$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").Value = "'160.73"
$ws.Range("E26").Value = "  +2.70%  "

$ws.Range("D27").Value = "'20.41"
$ws.Range("E27").Value = "  -4.62%  "

$ws.Range("D28").Value = "1.998.74"
$ws.Range("E28").Value = "  -1.73%  "

$ws.Range("D29").Value = "'2.331"
$ws.Range("E29").Value = "  -2.98%  "

$ws.Range("D30").Value = "'122.75"
$ws.Range("E30").Value = "  -2.62%  "

$ws.Range("D31").Value = "'0.1079"
$ws.Range("E31").Value = "  -1.15%  "

$ws.Range("E32").Value = "  -5.56%  "

$ws.Range("D33").Value = "'3.672"
$ws.Range("E33").Value = "  +0.36%  "

$ws.Range("D34").Value = "'5.536"
$ws.Range("E34").Value = "  -3.88%  "

$ws.Range("D35").Value = "'0.07226"
$ws.Range("E35").Value = "  +2.31%  "

$ws.Range("D36").Value = "'12.21"
$ws.Range("E36").Value = "  +8.10%  "

$ws.Range("D37").Value = "'0.02320"
$ws.Range("E37").Value = "  -1.23%  "

$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'8.801"
$ws.Range("E38").Value = "  -0.44%  "

$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2140"
$ws.Range("E39").Value = "  -4.00%  "

$ws.Range("D40").Value = "'5.057"
$ws.Range("E40").Value = "  -4.24%  "

$ws.Range("D41").Value = "'0.6141"
$ws.Range("E41").Value = "  -2.90%  "

$ws.Range("D42").Value = "'1.161"
$ws.Range("E42").Value = "  -1.81%  "

$ws.Range("D43").Value = "'1.355"
$ws.Range("E43").Value = "  -3.09%  "

$ws.Range("D44").Value = "'13.27"
$ws.Range("E44").Value = "  -1.79%  "

$ws.Range("D45").Value = "'3.772"
$ws.Range("E45").Value = "  +1.12%  "

$ws.Range("D46").Value = "'0.5969"
$ws.Range("E46").Value = "  +0.63%  "

$ws.Range("D47").Value = "'127.77"
$ws.Range("E47").Value = "  +2.03%  "

$ws.Range("D48").Value = "'1.231"
$ws.Range("E48").Value = "  +3.64%  "

$ws.Range("E49").Value = "  -3.74%  "

$ws.Range("D50").Value = "'0.06747"
$ws.Range("E50").Value = "  -2.28%  "

$ws.Range("D51").Value = "'72.96"
$ws.Range("E51").Value = "  -1.52%  "
